# Generate Report for Handoff
# Updates the localization-status workbook to reflect that the
# ca6f5501-68b9-41d7-94e5-4269e63add89 file is now "Ready for handoff",
# refreshes its timestamps, and records a handback version-mismatch error.

$wb = $excel.ActiveWorkbook

$status = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c370ad133225b882e04c5c4eb35b675ead8ba057/e2e/ca6f5501-68b9-41d7-94e5-4269e63add89.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/515f2d4b2689b96a50620a884941a82f8ef69ba8/e2e/ca6f5501-68b9-41d7-94e5-4269e63add89.md."

# --- Overview sheet: row 3 is the ca6f5501...md file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value2 = $status
$wsOverview.Range("F3").Value2 = $status
$wsOverview.Range("G3").Value2 = "2016-08-30 04:48:18"

# --- zh-cn sheet: row 3 is the ca6f5501...md file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value2 = $status
$wsZhCn.Range("H3").Value2 = "2016-08-30 04:48:13"
$wsZhCn.Range("P3").Value2 = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.15

# --- de-de sheet: row 3 is the ca6f5501...md file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value2 = $status
$wsDeDe.Range("H3").Value2 = "2016-08-30 04:48:18"
$wsDeDe.Range("P3").Value2 = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.15
